# Auto-generated edit script: updates crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.575.95'
$ws.Range("E2").Value = '  -0.81%  '
$ws.Range("D3").Value = '2.527.27'
$ws.Range("E3").Value = '  -1.83%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '308.20'
$ws.Range("E5").Value = '  -2.31%  '
$ws.Range("D6").Value = '99.82'
$ws.Range("E6").Value = '  +0.15%  '
$ws.Range("E7").Value = '  -1.59%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = '0.520'
$ws.Range("E9").Value = '  -2.95%  '
$ws.Range("D10").Value = '35.58'
$ws.Range("D11").Value = '0.0801'
$ws.Range("E11").Value = '  -1.42%  '
$ws.Range("D12").Value = '7.36'
$ws.Range("E12").Value = '  -2.19%  '
$ws.Range("E13").Value = '  +0.20%  '
$ws.Range("D14").Value = '2.916.35'
$ws.Range("E14").Value = '  -1.90%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '2.556.92'
$ws.Range("E15").Value = '  -3.20%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = '15.26'
$ws.Range("E16").Value = '  -3.03%  '
$ws.Range("D17").Value = '0.807'
$ws.Range("E17").Value = '  -4.43%  '
$ws.Range("D18").Value = '42.576.09'
$ws.Range("E19").Value = '  -2.22%  '
$ws.Range("D20").Value = '0.0₃0946'
$ws.Range("E20").Value = '  -2.30%  '
$ws.Range("D21").Value = '12.12'
$ws.Range("E21").Value = '  -4.13%  '
$ws.Range("D22").Value = '69.40'
$ws.Range("E22").Value = '  -0.22%  '
$ws.Range("D23").Value = '242.59'
$ws.Range("E23").Value = '  -2.94%  '
$ws.Range("D24").Value = '2.87'
$ws.Range("E24").Value = '  -3.41%  '
$ws.Range("D25").Value = '2.02'
$ws.Range("E25").Value = '  -2.79%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").Value = '25.45'
$ws.Range("E27").Value = '  -6.08%  '
$ws.Range("E28").Value = '  -2.41%  '
$ws.Range("D29").Value = '10.09'
$ws.Range("E29").Value = '  -2.05%  '
$ws.Range("D30").Value = '38.26'
$ws.Range("E30").Value = '  -5.75%  '
$ws.Range("D31").Value = '157.84'
$ws.Range("E31").Value = '  -0.30%  '
$ws.Range("D32").Value = '5.71'
$ws.Range("E32").Value = '  -2.08%  '
$ws.Range("E33").Value = '  +10.00%  '
$ws.Range("E34").Value = '  -1.62%  '
$ws.Range("D35").Value = '0.0780'
$ws.Range("E35").Value = '  -2.76%  '
$ws.Range("D37").Value = '17.88'
$ws.Range("E37").Value = '  -5.02%  '
$ws.Range("E38").Value = '  -7.70%  '
$ws.Range("E39").Value = '  -2.05%  '
$ws.Range("E40").Value = '  -1.08%  '
$ws.Range("D41").Value = '4.21'
$ws.Range("E41").Value = '  +2.43%  '
$ws.Range("D42").Value = '22.30'
$ws.Range("E42").Value = '  -5.47%  '
$ws.Range("E44").Value = '  -1.81%  '
$ws.Range("D45").Value = '3.25'
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("D46").Value = '2.003.28'
$ws.Range("E46").Value = '  -0.03%  '
$ws.Range("D47").Value = '8.87'
$ws.Range("E47").Value = '  -0.59%  '
$ws.Range("D48").Value = '2.769.47'
$ws.Range("E48").Value = '  -1.89%  '
$ws.Range("B49").Value = 'BitcoinSV'
$ws.Range("C49").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D49").Value = '79.08'
$ws.Range("E49").Value = '  -3.28%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '0.188'
$ws.Range("E50").Value = '  -4.80%  '
$ws.Range("D51").Value = '71.85'
$ws.Range("E51").Value = '  -4.23%  '
